$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -60
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 37
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = 37.037037037037
$ws.Range("L16").Value = 32.142857142857
$ws.Range("M16").Value = -21.276595744680
$ws.Range("N16").Value = -70.866141732283
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = -21.153846153846
$ws.Range("L17").Value = 10.810810810810
$ws.Range("M17").Value = 28.125
$ws.Range("N17").Value = 115.789473684211
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 160
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = 57.894736842105
$ws.Range("L18").Value = 36.363636363636
$ws.Range("M18").Value = -31.818181818181
$ws.Range("N18").Value = -87.755102040816
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = -24.390243902439
$ws.Range("I19").Value = 72
$ws.Range("J19").Value = 93
$ws.Range("K19").Value = -22.580645161290
$ws.Range("L19").Value = 26.315789473684
$ws.Range("M19").Value = 67.441860465116
$ws.Range("N19").Value = -1.369863013698
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 233.333333333333
$ws.Range("F20").Value = 36
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 38.461538461538
$ws.Range("I20").Value = 64
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = 39.130434782608
$ws.Range("L20").Value = 106.451612903226
$ws.Range("M20").Value = 120.689655172414
$ws.Range("N20").Value = -72.649572649572
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -32.258064516129
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = 0.787401574803
$ws.Range("I21").Value = 246
$ws.Range("J21").Value = 240
$ws.Range("K21").Value = 2.5
$ws.Range("L21").Value = 37.430167597765
$ws.Range("M21").Value = 24.242424242424
$ws.Range("N21").Value = -65.155807365439
$ws.Range("H22").Value = -100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 40
$ws.Range("K23").Value = 38.461538461538
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 80
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 74.025974025974
$ws.Range("I24").Value = 225
$ws.Range("J24").Value = 140
$ws.Range("K24").Value = 60.714285714285
$ws.Range("L24").Value = 42.405063291139
$ws.Range("M24").Value = 73.076923076923
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 75
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 31.428571428571
$ws.Range("I25").Value = 76
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = 16.923076923076
$ws.Range("L25").Value = 68.888888888888
$ws.Range("M25").Value = 8.571428571428
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 50
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 125
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = 200
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -25
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = 200
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0

# --- Cells whose underlying type/style also changes ---
# Text(General, s=14) -> Numeric: reuse style 15 (#,##0) from C16, or style 16 (pct) from E16
# Numeric -> Text(General, s=14): reuse style 14 from A16 (apostrophe forces text entry)
$ws.Range("D18").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E18").Value = -50
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F22").Value = "'0"
$ws.Range("A16").Copy()
$ws.Range("F22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D23").Value = "'0"
$ws.Range("A16").Copy()
$ws.Range("D23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E23").Value = "'***.*"
$ws.Range("A16").Copy()
$ws.Range("E23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C26").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C28").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C29").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false
